$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price / Volume(1h) / Hora columns hold numeric-looking text
# ("307.00", "2.01%", "3", ...). Writing those strings straight into a
# General-formatted cell makes Excel auto-convert them to real numbers
# (dropping significant trailing zeros, turning "2.01%" into 0.0201,
# etc.), same as typing them in the Excel UI. Mark those cells as Text
# first so the values land verbatim, exactly like the source data.
$textCells = @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "E24", "G24", "E25", "G25", "D26", "E26", "G26", "E27", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46", "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50", "G50", "D51", "E51", "G51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "307.00"
$ws.Range("E2").Value = "2.01%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "36.33"
$ws.Range("E3").Value = "-4.96%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "5.039"
$ws.Range("E4").Value = "0.98%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.07837"
$ws.Range("E5").Value = "1.67%"
$ws.Range("G5").Value = "3"
$ws.Range("D6").Value = "2.121"
$ws.Range("E6").Value = "-3.17%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "7.928"
$ws.Range("E7").Value = "-0.49%"
$ws.Range("G7").Value = "3"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9210"
$ws.Range("E8").Value = "0.47%"
$ws.Range("G8").Value = "3"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09608"
$ws.Range("E9").Value = "6.78%"
$ws.Range("G9").Value = "3"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1887"
$ws.Range("E10").Value = "5.63%"
$ws.Range("G10").Value = "3"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08736"
$ws.Range("E11").Value = "3.46%"
$ws.Range("G11").Value = "3"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03500"
$ws.Range("E12").Value = "-1.21%"
$ws.Range("G12").Value = "3"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09931"
$ws.Range("E13").Value = "-0.15%"
$ws.Range("G13").Value = "3"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001428"
$ws.Range("E14").Value = "-3.14%"
$ws.Range("G14").Value = "3"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005710"
$ws.Range("E15").Value = "1.08%"
$ws.Range("G15").Value = "3"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.460"
$ws.Range("E16").Value = "-0.55%"
$ws.Range("G16").Value = "3"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.063"
$ws.Range("E17").Value = "1.74%"
$ws.Range("G17").Value = "3"
$ws.Range("D18").Value = "2.409"
$ws.Range("E18").Value = "8.33%"
$ws.Range("G18").Value = "3"
$ws.Range("G19").Value = "3"
$ws.Range("D20").Value = "0.1337"
$ws.Range("E20").Value = "2.72%"
$ws.Range("G20").Value = "3"
$ws.Range("D21").Value = "4.772"
$ws.Range("E21").Value = "4.67%"
$ws.Range("G21").Value = "3"
$ws.Range("D22").Value = "0.2291"
$ws.Range("E22").Value = "2.22%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04595"
$ws.Range("E23").Value = "-1.43%"
$ws.Range("G23").Value = "3"
$ws.Range("E24").Value = "15.04%"
$ws.Range("G24").Value = "3"
$ws.Range("E25").Value = "-0.04%"
$ws.Range("G25").Value = "3"
$ws.Range("D26").Value = "0.0001399"
$ws.Range("E26").Value = "7.39%"
$ws.Range("G26").Value = "3"
$ws.Range("E27").Value = "-42.74%"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("D39").Value = "0.01835"
$ws.Range("E39").Value = "5.71%"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.04794"
$ws.Range("E40").Value = "2.65%"
$ws.Range("G40").Value = "3"
$ws.Range("D41").Value = "0.007469"
$ws.Range("E41").Value = "-4.91%"
$ws.Range("G41").Value = "3"
$ws.Range("D42").Value = "0.1402"
$ws.Range("E42").Value = "1.12%"
$ws.Range("G42").Value = "3"
$ws.Range("D43").Value = "0.007740"
$ws.Range("E43").Value = "0.74%"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.002236"
$ws.Range("E44").Value = "-2.57%"
$ws.Range("G44").Value = "3"
$ws.Range("D45").Value = "0.01039"
$ws.Range("E45").Value = "8.75%"
$ws.Range("G45").Value = "3"
$ws.Range("D46").Value = "0.00006165"
$ws.Range("E46").Value = "1.81%"
$ws.Range("G46").Value = "3"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.28%"
$ws.Range("G47").Value = "3"
$ws.Range("D48").Value = "0.0005803"
$ws.Range("E48").Value = "0.04%"
$ws.Range("G48").Value = "3"
$ws.Range("D49").Value = "42.78"
$ws.Range("E49").Value = "402.08%"
$ws.Range("G49").Value = "3"
$ws.Range("D50").Value = "0.002001"
$ws.Range("E50").Value = "-25.85%"
$ws.Range("G50").Value = "3"
$ws.Range("D51").Value = "0.00002099"
$ws.Range("E51").Value = "-0.28%"
$ws.Range("G51").Value = "3"
